$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the next empty row right after the last row of data in column A.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# New scrape results for this run (Date, Price, Discount, Incredible).
$scrapeDate = "2026-02-07"
$price      = "33031200"
$discount   = "0"
$incredible = "0"

# Values must land in the sheet as plain text (shared strings), exactly
# like every other row, instead of being auto-coerced by Excel into a
# date serial / numeric value. Using a leading apostrophe forces text
# entry; we then reset the cell style back to the workbook's default
# "Normal" style so no extra number-format styling is left behind.
$rng = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 4))

$ws.Cells.Item($newRow, 1).Value = "'" + $scrapeDate
$ws.Cells.Item($newRow, 2).Value = "'" + $price
$ws.Cells.Item($newRow, 3).Value = "'" + $discount
$ws.Cells.Item($newRow, 4).Value = "'" + $incredible

$rng.Style = "Normal"
